$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.381182789802551
$ws.Range("B1").Value = 2.56744909286499
$ws.Range("C1").Value = 6.541423797607422
$ws.Range("D1").Value = 2.403305768966675
$ws.Range("E1").Value = 1.214823246002197
